$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 436, pushing the existing
# rows 436-440 down to 439-443 (their original content is preserved
# automatically by the insert).
$ws.Rows("436:438").Insert()

# Row 436: new "Especial" quality entry
$ws.Cells.Item(436,1).Value = 4
$ws.Cells.Item(436,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(436,3).Value = "Los Lagos"
$ws.Cells.Item(436,4).Value = 44890
$ws.Cells.Item(436,5).Value = 10
$ws.Cells.Item(436,6).Value = "Fruta"
$ws.Cells.Item(436,7).Value = 100102
$ws.Cells.Item(436,8).Value = "Cítricos"
$ws.Cells.Item(436,9).Value = 100102006
$ws.Cells.Item(436,10).Value = "Pomelo"
$ws.Cells.Item(436,11).Value = "Start Ruby"
$ws.Cells.Item(436,12).Value = "Especial"
$ws.Cells.Item(436,13).Value = 50
$ws.Cells.Item(436,14).Value = 16000
$ws.Cells.Item(436,15).Value = 16000
$ws.Cells.Item(436,16).Value = 16000
$ws.Cells.Item(436,17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(436,18).Value = "Región de O'Higgins"
$ws.Cells.Item(436,19).Value = 1143
$ws.Cells.Item(436,20).Value = 14

# Row 437: new "Primera" quality entry
$ws.Cells.Item(437,1).Value = 4
$ws.Cells.Item(437,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(437,3).Value = "Los Lagos"
$ws.Cells.Item(437,4).Value = 44890
$ws.Cells.Item(437,5).Value = 10
$ws.Cells.Item(437,6).Value = "Fruta"
$ws.Cells.Item(437,7).Value = 100102
$ws.Cells.Item(437,8).Value = "Cítricos"
$ws.Cells.Item(437,9).Value = 100102006
$ws.Cells.Item(437,10).Value = "Pomelo"
$ws.Cells.Item(437,11).Value = "Start Ruby"
$ws.Cells.Item(437,12).Value = "Primera"
$ws.Cells.Item(437,13).Value = 50
$ws.Cells.Item(437,14).Value = 14000
$ws.Cells.Item(437,15).Value = 14000
$ws.Cells.Item(437,16).Value = 14000
$ws.Cells.Item(437,17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(437,18).Value = "Región de O'Higgins"
$ws.Cells.Item(437,19).Value = 1000
$ws.Cells.Item(437,20).Value = 14

# Row 438: new "Segunda" quality entry
$ws.Cells.Item(438,1).Value = 4
$ws.Cells.Item(438,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(438,3).Value = "Los Lagos"
$ws.Cells.Item(438,4).Value = 44890
$ws.Cells.Item(438,5).Value = 10
$ws.Cells.Item(438,6).Value = "Fruta"
$ws.Cells.Item(438,7).Value = 100102
$ws.Cells.Item(438,8).Value = "Cítricos"
$ws.Cells.Item(438,9).Value = 100102006
$ws.Cells.Item(438,10).Value = "Pomelo"
$ws.Cells.Item(438,11).Value = "Start Ruby"
$ws.Cells.Item(438,12).Value = "Segunda"
$ws.Cells.Item(438,13).Value = 50
$ws.Cells.Item(438,14).Value = 12000
$ws.Cells.Item(438,15).Value = 12000
$ws.Cells.Item(438,16).Value = 12000
$ws.Cells.Item(438,17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(438,18).Value = "Región de O'Higgins"
$ws.Cells.Item(438,19).Value = 857
$ws.Cells.Item(438,20).Value = 14
